$wb = $excel.ActiveWorkbook

# ---- Metadata sheet updates ----
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B5").Value = "T-CABS ValueSet Ventilation Mode"
$meta.Range("B8").Value = "2025-11-19T11:55:29+01:00"
$meta.Range("B10").Value = "BIH-CEI (https://www.bihealth.org/)"
$meta.Range("B11").Value = "Germany"
$meta.Range("B12").Value = "This ValueSet contains codes to represent the different ventilation modes"

# ---- Include #0 (concept table) sheet updates ----
$inc = $wb.Worksheets.Item("Include #0")

# Replace the old trailing concept row (151796 / MDC_PRESS_AWAY_CTS_POS)
# with the next entry in the new sequence, then append the remaining
# new concept rows, the blank separator row, and the System URI row.
$inc.Range("A8").Value = "475161"
$inc.Range("B8").Value = "MDC_VENT_MODE_ISO_CSV_6ACAP"

$inc.Range("A9").Value = "475140"
$inc.Range("B9").Value = "MDC_VENT_MODE_ISO_AC_VC_6ACAPL"

$inc.Range("A10").Value = "475136"
$inc.Range("B10").Value = "MDC_VENT_MODE_ISO_SIMV_VC_8PS_6ACAPL"

$inc.Range("A11").Value = "475190"
$inc.Range("B11").Value = "MDC_VENT_MODE_ISO_HIFLOW_3NIV"

$inc.Range("A12").Value = ""
$inc.Range("B12").Value = ""

$inc.Range("A13").Value = "System URI"
$inc.Range("B13").Value = "urn:iso:std:iso:11073:10101"
